$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Checkmark glyph already used elsewhere in the "Konfirmasi" column (U+2714 U+FE0F)
$checkmark = "✔️"

# Mark confirmation checkmarks for the Ustadz Hanif rows (27, 29, 31, 33)
$ws.Range("E27").Value = $checkmark
$ws.Range("E29").Value = $checkmark
$ws.Range("E31").Value = $checkmark
$ws.Range("E33").Value = $checkmark

# Reflect where the user ended up scrolled to / selecting when they made the edit.
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E27").Select()
